# Add a new worksheet "strategy_id-7160" at the end of the workbook,
# mirroring the structure/content of the existing "strategy_id-7106" sheet
# (bold/bordered header row of labels + numbers, plus one data row),
# matching the commit diff.

$wb = $excel.ActiveWorkbook

# The existing sheet whose layout/formatting the new sheet should mirror.
$template = $wb.Worksheets.Item("strategy_id-7106")

# Insert the new sheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "strategy_id-7160"

# ---- Row 1: header labels (A1:I1) and index numbers (J1:AS1) -------------
$headerLabels = @("subsector","variable","variable_trajectory_group","normalize_group","trajgroup_no_vary_q","uniform_scaling_q","variable_trajectory_group_trajectory_type","max_35","min_35")

for ($i = 0; $i -lt $headerLabels.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerLabels[$i]
}

for ($n = 0; $n -le 35; $n++) {
    $ws.Cells.Item(1, 10 + $n).Value = $n
}

# Mirror the bold / centered / bordered header formatting from the template
# sheet instead of re-deriving it property-by-property (keeps styles.xml
# free of extra, unused cellXfs entries).
$template.Range("A1:AS1").Copy()
$ws.Range("A1:AS1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 2: data row ------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "General"
$ws.Cells.Item(2, 2).Value = "frac_gnrl_eating_red_meat"
$ws.Cells.Item(2, 3).Value = 13

for ($c = 8; $c -le 22; $c++) {
    $ws.Cells.Item(2, $c).Value = 1
}

$decayValues = @(
    0.9869565217391305,
    0.9739130434782609,
    0.9608695652173913,
    0.9478260869565217,
    0.9347826086956522,
    0.9217391304347825,
    0.908695652173913,
    0.8956521739130434,
    0.8826086956521739,
    0.8695652173913044,
    0.8565217391304347,
    0.8434782608695652,
    0.8304347826086956,
    0.817391304347826,
    0.8043478260869565,
    0.7913043478260869,
    0.7782608695652173,
    0.7652173913043477,
    0.7521739130434782,
    0.7391304347826086,
    0.7260869565217392,
    0.7130434782608696,
    0.7
)

for ($i = 0; $i -lt $decayValues.Count; $i++) {
    $ws.Cells.Item(2, 23 + $i).Value = $decayValues[$i]
}

# Restore the originally-active first sheet as the selected tab so adding
# this sheet doesn't perturb the workbook's prior tab-selection state.
$wb.Worksheets.Item(1).Activate()
